$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(4).Insert()

$ws.Range("D1").Value = "participantsCount"
$ws.Range("D2").Value = 11
$ws.Range("D3").Value = 11
$ws.Range("D4").Value = 11

$ws.Range("E2").Value = "Veverky"
$ws.Range("E3").Value = "Bobři"
$ws.Range("E4").Value = "Křečci"
